# Commit: "Update after the Ondra Kucera merge request"
#
# 1) Drop slides 2-4, keeping only the title slide.
# 2) Merge the title's two text runs ("Java 1 - Lekce " + "03") into one run.
# 3) Remove now-unused "Nadpis a obsah" layout (only used by the deleted slides).
# 4) Strip the old sponsor/partner logo pictures from the remaining layouts.

$p = $ppt.ActivePresentation

# --- 1) Remove the extra slides (keep slide 1 only) -------------------------
for ($i = $p.Slides.Count; $i -ge 2; $i--) {
    $p.Slides.Item($i).Delete()
}

# --- 2) Merge the title run on slide 1 --------------------------------------
$slide1 = $p.Slides.Item(1)
$title = $slide1.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$dash = [char]0x2013
$mergedText = "Java 1 " + $dash + " Lekce 03"
$wholeRange = $titleRange.Characters(1, $titleRange.Length)
$wholeRange.Text = $mergedText

# --- 3) Delete the now-orphaned "Nadpis a obsah" slide layout --------------
$master = $p.SlideMaster
for ($i = $master.CustomLayouts.Count; $i -ge 1; $i--) {
    $layout = $master.CustomLayouts.Item($i)
    if ($layout.Name -eq "Nadpis a obsah") {
        $layout.Delete()
    }
}

# --- 4) Remove the obsolete logo pictures from the remaining layouts -------
$picNamesToRemove = @("Shape 12", "Shape 16", "Shape 20", "Picture 2", "Shape 30", "Shape 39", "Shape 54")
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = $layout.Shapes.Count; $j -ge 1; $j--) {
        $shape = $layout.Shapes.Item($j)
        if ($picNamesToRemove -contains $shape.Name) {
            $shape.Delete()
        }
    }
}

Write-Output ("Final slide count: " + $p.Slides.Count)
Write-Output ("Final layout count: " + $master.CustomLayouts.Count)
Write-Output ("Final title text: " + $title.TextFrame.TextRange.Text)
